# Weekly update: insert a new week's worth of Tomate price records
# (Comercializadora del Agro de Limari) ahead of the existing history,
# pushing the remaining rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new blank rows starting at row 1057; this shifts the existing
# rows 1057:1096 down to 1062:1101 and updates the sheet dimension
# automatically.
$ws.Rows("1057:1061").Insert()

# Static column values shared by every record in this data block.
$mercadoId  = 2
$mercado    = "Comercializadora del Agro de Limarí"
$region     = "Coquimbo"
$codreg     = 4
$categoriaId = 100112020
$categoria  = "Tomate"
$unidad     = "`$/bandeja 18 kilos"
$origen     = "Provincia de Limarí"
$kgUnidades = 18
$clasificacion = "Hortaliza"

# New-week data: Fecha (serial), Variedad, Calidad, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg
$newRows = @(
    @{ Row=1057; D=45267; H="Larga vida"; I="Primera"; J=2000; K=11000; L=12000; M=11500; P=639 },
    @{ Row=1058; D=45267; H="Larga vida"; I="Segunda"; J=1800; K=8000;  L=9000;  M=8500;  P=472 },
    @{ Row=1059; D=45267; H="Larga vida"; I="Tercera"; J=1400; K=5000;  L=6000;  M=5500;  P=306 },
    @{ Row=1060; D=45267; H="Semiduro";   I="Primera"; J=1800; K=8000;  L=9000;  M=8500;  P=472 },
    @{ Row=1061; D=45267; H="Semiduro";   I="Segunda"; J=1000; K=5000;  L=6000;  M=5500;  P=306 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value2  = $mercadoId
    $ws.Cells.Item($row, 2).Value2  = $mercado
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = $r.D
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $categoriaId
    $ws.Cells.Item($row, 7).Value2  = $categoria
    $ws.Cells.Item($row, 8).Value2  = $r.H
    $ws.Cells.Item($row, 9).Value2  = $r.I
    $ws.Cells.Item($row, 10).Value2 = $r.J
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = $r.M
    $ws.Cells.Item($row, 14).Value2 = $unidad
    $ws.Cells.Item($row, 15).Value2 = $origen
    $ws.Cells.Item($row, 16).Value2 = $r.P
    $ws.Cells.Item($row, 17).Value2 = $kgUnidades
    $ws.Cells.Item($row, 18).Value2 = $clasificacion
}
